# Applies three edits to the document:
#  1. In the "Documento:" table row, merge the two runs
#     "Plan de acción del Proyecto " + "Capstone" into one run/text,
#     which also drops the spell-check proofErr markers around "Capstone".
#  2. In the document title paragraph, merge the two runs
#     "Plan de acción del proyecto " + "Capstone" into one run/text,
#     same proofErr cleanup.
#  3. Fill in the previously empty "Numero de equipo" value cell with "14".

$d = $word.ActiveDocument

# --- 1 & 2: merge the split "... Capstone" runs via Find & Replace ---
# Word's Find/Replace rewrites the paragraph's run content as a single
# run of the replacement text, which naturally removes the now-stale
# proofErr spell-check bookmarks that bracketed "Capstone".

$d.Content.Find.Execute(
    "Plan de acción del Proyecto Capstone", $true, $false, $false, $false,
    $false, $true, 1, $false, "Plan de acción del Proyecto Capstone", 2)

$d.Content.Find.Execute(
    "Plan de acción del proyecto Capstone", $true, $false, $false, $false,
    $false, $true, 1, $false, "Plan de acción del proyecto Capstone", 2)

# --- 3: set the team-number cell text to "14" ---
$table = $d.Tables.Item(2)
$cell = $table.Cell(3, 2)
$cell.Range.Text = "14"
